$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Part 1: ten pairs of adjacent rows had their match data (everything
# except id/Div/Div Original Name/Date) swapped between the two rows.
# ------------------------------------------------------------------
# --- swap rows 19 and 20 ---
$ws.Range("B19").Value = 6100758
$ws.Range("B20").Value = 6100756
$ws.Range("F19").Value = 'Umea FC'
$ws.Range("F20").Value = 'Sollentuna United FF'
$ws.Range("G19").Value = 'Pite IF'
$ws.Range("G20").Value = 'IFK Stocksund'
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 3
$ws.Range("I19").Value = 0
$ws.Range("I20").Value = 1
$ws.Range("J19").Value = 'H'
$ws.Range("J20").Value = 'H'
$ws.Range("K19").Value = 1.909
$ws.Range("K20").Value = 1.666
$ws.Range("L19").Value = 3.75
$ws.Range("L20").Value = 3.75
$ws.Range("M19").Value = 3.1
$ws.Range("M20").Value = 4
$ws.Range("N19").Value = 2.25
$ws.Range("N20").Value = 1.7
$ws.Range("O19").Value = 3.6
$ws.Range("O20").Value = 3.75
$ws.Range("P19").Value = 2.625
$ws.Range("P20").Value = 3.8
$ws.Range("Q19").Value = -0.25
$ws.Range("Q20").Value = -0.75
$ws.Range("R19").Value = 2.025
$ws.Range("R20").Value = 1.95
$ws.Range("S19").Value = 1.775
$ws.Range("S20").Value = 1.85
$ws.Range("T19").Value = 3
$ws.Range("T20").Value = 3.25
$ws.Range("U19").Value = 1.875
$ws.Range("U20").Value = 1.85
$ws.Range("V19").Value = 1.925
$ws.Range("V20").Value = 1.95
$ws.Range("W19").Value = 1.25
$ws.Range("W20").Value = 0.7
$ws.Range("X19").Value = -1
$ws.Range("X20").Value = -1
$ws.Range("Y19").Value = -1
$ws.Range("Y20").Value = -1
$ws.Range("Z19").Value = 1.025
$ws.Range("Z20").Value = 0.95
$ws.Range("AA19").Value = -1
$ws.Range("AA20").Value = -1
$ws.Range("AB19").Value = -1
$ws.Range("AB20").Value = 0.8500000000000001
$ws.Range("AC19").Value = 0.925
$ws.Range("AC20").Value = -1

# --- swap rows 32 and 33 ---
$ws.Range("B32").Value = 6800952
$ws.Range("B33").Value = 6100763
$ws.Range("F32").Value = 'IF Karlstad Fotboll'
$ws.Range("F33").Value = 'rebro Syrianska IF'
$ws.Range("G32").Value = 'Dalkurd FF'
$ws.Range("G33").Value = 'Sollentuna United FF'
$ws.Range("H32").Value = 2
$ws.Range("H33").Value = 1
$ws.Range("I32").Value = 1
$ws.Range("I33").Value = 1
$ws.Range("J32").Value = 'H'
$ws.Range("J33").Value = 'D'
$ws.Range("K32").Value = 3.6
$ws.Range("K33").Value = 2.2
$ws.Range("L32").Value = 3.6
$ws.Range("L33").Value = 3.5
$ws.Range("M32").Value = 1.8
$ws.Range("M33").Value = 2.7
$ws.Range("N32").Value = 3.3
$ws.Range("N33").Value = 2.75
$ws.Range("O32").Value = 3.6
$ws.Range("O33").Value = 3.5
$ws.Range("P32").Value = 1.909
$ws.Range("P33").Value = 2.2
$ws.Range("Q32").Value = 0.5
$ws.Range("Q33").Value = 0.25
$ws.Range("R32").Value = 1.85
$ws.Range("R33").Value = 1.8
$ws.Range("S32").Value = 1.95
$ws.Range("S33").Value = 2
$ws.Range("T32").Value = 2.75
$ws.Range("T33").Value = 2.75
$ws.Range("U32").Value = 1.8
$ws.Range("U33").Value = 1.85
$ws.Range("V32").Value = 2
$ws.Range("V33").Value = 1.95
$ws.Range("W32").Value = 2.3
$ws.Range("W33").Value = -1
$ws.Range("X32").Value = -1
$ws.Range("X33").Value = 2.5
$ws.Range("Y32").Value = -1
$ws.Range("Y33").Value = -1
$ws.Range("Z32").Value = 0.8500000000000001
$ws.Range("Z33").Value = 0.4
$ws.Range("AA32").Value = -1
$ws.Range("AA33").Value = -0.5
$ws.Range("AB32").Value = 0.4
$ws.Range("AB33").Value = -1
$ws.Range("AC32").Value = -0.5
$ws.Range("AC33").Value = 0.95

# --- swap rows 50 and 51 ---
$ws.Range("B50").Value = 7048209
$ws.Range("B51").Value = 6100778
$ws.Range("F50").Value = 'Hammarby TFF'
$ws.Range("F51").Value = 'Nordic United FC'
$ws.Range("G50").Value = 'IF Karlstad Fotboll'
$ws.Range("G51").Value = 'FC Stockholm Internazionale'
$ws.Range("H50").Value = 1
$ws.Range("H51").Value = 4
$ws.Range("I50").Value = 1
$ws.Range("I51").Value = 2
$ws.Range("J50").Value = 'D'
$ws.Range("J51").Value = 'H'
$ws.Range("K50").Value = 2.25
$ws.Range("K51").Value = 1.75
$ws.Range("L50").Value = 3.5
$ws.Range("L51").Value = 3.6
$ws.Range("M50").Value = 2.6
$ws.Range("M51").Value = 3.75
$ws.Range("N50").Value = 2.5
$ws.Range("N51").Value = 1.75
$ws.Range("O50").Value = 3.3
$ws.Range("O51").Value = 3.5
$ws.Range("P50").Value = 2.4
$ws.Range("P51").Value = 4
$ws.Range("Q50").Value = 0
$ws.Range("Q51").Value = -0.75
$ws.Range("R50").Value = 1.975
$ws.Range("R51").Value = 2
$ws.Range("S50").Value = 1.825
$ws.Range("S51").Value = 1.8
$ws.Range("T50").Value = 2.5
$ws.Range("T51").Value = 2.75
$ws.Range("U50").Value = 1.9
$ws.Range("U51").Value = 1.825
$ws.Range("V50").Value = 1.9
$ws.Range("V51").Value = 1.975
$ws.Range("W50").Value = -1
$ws.Range("W51").Value = 0.75
$ws.Range("X50").Value = 2.3
$ws.Range("X51").Value = -1
$ws.Range("Y50").Value = -1
$ws.Range("Y51").Value = -1
$ws.Range("Z50").Value = 0
$ws.Range("Z51").Value = 1
$ws.Range("AA50").Value = 0
$ws.Range("AA51").Value = -1
$ws.Range("AB50").Value = -1
$ws.Range("AB51").Value = 0.825
$ws.Range("AC50").Value = 0.8999999999999999
$ws.Range("AC51").Value = -1

# --- swap rows 56 and 57 ---
$ws.Range("B56").Value = 6100513
$ws.Range("B57").Value = 6100775
$ws.Range("F56").Value = 'Umea FC'
$ws.Range("F57").Value = 'IF Sylvia'
$ws.Range("G56").Value = 'Bodens BK FF'
$ws.Range("G57").Value = 'Sandvikens IF'
$ws.Range("H56").Value = 2
$ws.Range("H57").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("I57").Value = 3
$ws.Range("J56").Value = 'H'
$ws.Range("J57").Value = 'A'
$ws.Range("K56").Value = 1.615
$ws.Range("K57").Value = 5.75
$ws.Range("L56").Value = 3.75
$ws.Range("L57").Value = 4.2
$ws.Range("M56").Value = 4.333
$ws.Range("M57").Value = 1.4
$ws.Range("N56").Value = 1.727
$ws.Range("N57").Value = 6.5
$ws.Range("O56").Value = 3.6
$ws.Range("O57").Value = 4.75
$ws.Range("P56").Value = 4
$ws.Range("P57").Value = 1.333
$ws.Range("Q56").Value = -0.75
$ws.Range("Q57").Value = 1.5
$ws.Range("R56").Value = 2
$ws.Range("R57").Value = 1.925
$ws.Range("S56").Value = 1.8
$ws.Range("S57").Value = 1.875
$ws.Range("T56").Value = 2.75
$ws.Range("T57").Value = 3.5
$ws.Range("U56").Value = 1.85
$ws.Range("U57").Value = 1.875
$ws.Range("V56").Value = 1.95
$ws.Range("V57").Value = 1.925
$ws.Range("W56").Value = 0.7270000000000001
$ws.Range("W57").Value = -1
$ws.Range("X56").Value = -1
$ws.Range("X57").Value = -1
$ws.Range("Y56").Value = -1
$ws.Range("Y57").Value = 0.333
$ws.Range("Z56").Value = 1
$ws.Range("Z57").Value = -1
$ws.Range("AA56").Value = -1
$ws.Range("AA57").Value = 0.875
$ws.Range("AB56").Value = -1
$ws.Range("AB57").Value = -1
$ws.Range("AC56").Value = 0.95
$ws.Range("AC57").Value = 0.925

# --- swap rows 61 and 62 ---
$ws.Range("B61").Value = 6097795
$ws.Range("B62").Value = 6100779
$ws.Range("F61").Value = 'Dalkurd FF'
$ws.Range("F62").Value = 'IFK Stocksund'
$ws.Range("G61").Value = 'Tby FK'
$ws.Range("G62").Value = 'Motala AIF FK'
$ws.Range("H61").Value = 3
$ws.Range("H62").Value = 2
$ws.Range("I61").Value = 2
$ws.Range("I62").Value = 2
$ws.Range("J61").Value = 'H'
$ws.Range("J62").Value = 'D'
$ws.Range("K61").Value = 1.444
$ws.Range("K62").Value = 2.2
$ws.Range("L61").Value = 4.2
$ws.Range("L62").Value = 3.6
$ws.Range("M61").Value = 6.5
$ws.Range("M62").Value = 2.7
$ws.Range("N61").Value = 1.4
$ws.Range("N62").Value = 2.05
$ws.Range("O61").Value = 4.333
$ws.Range("O62").Value = 3.8
$ws.Range("P61").Value = 7
$ws.Range("P62").Value = 2.9
$ws.Range("Q61").Value = -1.25
$ws.Range("Q62").Value = -0.25
$ws.Range("R61").Value = 1.825
$ws.Range("R62").Value = 1.85
$ws.Range("S61").Value = 1.975
$ws.Range("S62").Value = 1.95
$ws.Range("T61").Value = 3.25
$ws.Range("T62").Value = 3.25
$ws.Range("U61").Value = 1.9
$ws.Range("U62").Value = 1.95
$ws.Range("V61").Value = 1.9
$ws.Range("V62").Value = 1.85
$ws.Range("W61").Value = 0.3999999999999999
$ws.Range("W62").Value = -1
$ws.Range("X61").Value = -1
$ws.Range("X62").Value = 2.8
$ws.Range("Y61").Value = -1
$ws.Range("Y62").Value = -1
$ws.Range("Z61").Value = -0.5
$ws.Range("Z62").Value = -0.5
$ws.Range("AA61").Value = 0.4875
$ws.Range("AA62").Value = 0.475
$ws.Range("AB61").Value = 0.8999999999999999
$ws.Range("AB62").Value = 0.95
$ws.Range("AC61").Value = -1
$ws.Range("AC62").Value = -1

# --- swap rows 66 and 67 ---
$ws.Range("B66").Value = 7115552
$ws.Range("B67").Value = 6097796
$ws.Range("F66").Value = 'rebro Syrianska IF'
$ws.Range("F67").Value = 'Sollentuna United FF'
$ws.Range("G66").Value = 'IF Karlstad Fotboll'
$ws.Range("G67").Value = 'Dalkurd FF'
$ws.Range("H66").Value = 0
$ws.Range("H67").Value = 3
$ws.Range("I66").Value = 3
$ws.Range("I67").Value = 1
$ws.Range("J66").Value = 'A'
$ws.Range("J67").Value = 'H'
$ws.Range("K66").Value = 3.25
$ws.Range("K67").Value = 2.75
$ws.Range("L66").Value = 3.4
$ws.Range("L67").Value = 3.5
$ws.Range("M66").Value = 2
$ws.Range("M67").Value = 2.25
$ws.Range("N66").Value = 3.75
$ws.Range("N67").Value = 2.2
$ws.Range("O66").Value = 3.6
$ws.Range("O67").Value = 3.75
$ws.Range("P66").Value = 1.8
$ws.Range("P67").Value = 2.625
$ws.Range("Q66").Value = 0.5
$ws.Range("Q67").Value = -0.25
$ws.Range("R66").Value = 1.95
$ws.Range("R67").Value = 2
$ws.Range("S66").Value = 1.85
$ws.Range("S67").Value = 1.8
$ws.Range("T66").Value = 3
$ws.Range("T67").Value = 3
$ws.Range("U66").Value = 1.95
$ws.Range("U67").Value = 1.825
$ws.Range("V66").Value = 1.85
$ws.Range("V67").Value = 1.975
$ws.Range("W66").Value = -1
$ws.Range("W67").Value = 1.2
$ws.Range("X66").Value = -1
$ws.Range("X67").Value = -1
$ws.Range("Y66").Value = 0.8
$ws.Range("Y67").Value = -1
$ws.Range("Z66").Value = -1
$ws.Range("Z67").Value = 1
$ws.Range("AA66").Value = 0.8500000000000001
$ws.Range("AA67").Value = -1
$ws.Range("AB66").Value = 0
$ws.Range("AB67").Value = 0.825
$ws.Range("AC66").Value = 0
$ws.Range("AC67").Value = -1

# --- swap rows 82 and 84 ---
$ws.Range("B82").Value = 6100805
$ws.Range("B84").Value = 6100808
$ws.Range("F82").Value = 'Motala AIF FK'
$ws.Range("F84").Value = 'rebro Syrianska IF'
$ws.Range("G82").Value = 'Vasalunds IF'
$ws.Range("G84").Value = 'Nordic United FC'
$ws.Range("H82").Value = 1
$ws.Range("H84").Value = 4
$ws.Range("I82").Value = 2
$ws.Range("I84").Value = 4
$ws.Range("J82").Value = 'A'
$ws.Range("J84").Value = 'D'
$ws.Range("K82").Value = 3.1
$ws.Range("K84").Value = 3.4
$ws.Range("L82").Value = 3.4
$ws.Range("L84").Value = 3.8
$ws.Range("M82").Value = 2
$ws.Range("M84").Value = 1.8
$ws.Range("N82").Value = 3.4
$ws.Range("N84").Value = 6
$ws.Range("O82").Value = 3.4
$ws.Range("O84").Value = 4.75
$ws.Range("P82").Value = 1.909
$ws.Range("P84").Value = 1.363
$ws.Range("Q82").Value = 0.5
$ws.Range("Q84").Value = 1.25
$ws.Range("R82").Value = 1.8
$ws.Range("R84").Value = 1.95
$ws.Range("S82").Value = 2
$ws.Range("S84").Value = 1.85
$ws.Range("T82").Value = 2.75
$ws.Range("T84").Value = 3.25
$ws.Range("U82").Value = 1.95
$ws.Range("U84").Value = 1.975
$ws.Range("V82").Value = 1.85
$ws.Range("V84").Value = 1.825
$ws.Range("W82").Value = -1
$ws.Range("W84").Value = -1
$ws.Range("X82").Value = -1
$ws.Range("X84").Value = 3.75
$ws.Range("Y82").Value = 0.909
$ws.Range("Y84").Value = -1
$ws.Range("Z82").Value = -1
$ws.Range("Z84").Value = 0.95
$ws.Range("AA82").Value = 1
$ws.Range("AA84").Value = -1
$ws.Range("AB82").Value = 0.475
$ws.Range("AB84").Value = 0.9750000000000001
$ws.Range("AC82").Value = -0.5
$ws.Range("AC84").Value = -1

# --- swap rows 93 and 94 ---
$ws.Range("B93").Value = 6100812
$ws.Range("B94").Value = 6100813
$ws.Range("F93").Value = 'Tby FK'
$ws.Range("F94").Value = 'Vasalunds IF'
$ws.Range("G93").Value = 'Sandvikens IF'
$ws.Range("G94").Value = 'rebro Syrianska IF'
$ws.Range("H93").Value = 0
$ws.Range("H94").Value = 2
$ws.Range("I93").Value = 3
$ws.Range("I94").Value = 0
$ws.Range("J93").Value = 'A'
$ws.Range("J94").Value = 'H'
$ws.Range("K93").Value = 5.75
$ws.Range("K94").Value = 1.45
$ws.Range("L93").Value = 4.5
$ws.Range("L94").Value = 4
$ws.Range("M93").Value = 1.4
$ws.Range("M94").Value = 5.5
$ws.Range("N93").Value = 7
$ws.Range("N94").Value = 1.363
$ws.Range("O93").Value = 5
$ws.Range("O94").Value = 4.5
$ws.Range("P93").Value = 1.3
$ws.Range("P94").Value = 6.5
$ws.Range("Q93").Value = 1.25
$ws.Range("Q94").Value = -1.25
$ws.Range("R93").Value = 2.1
$ws.Range("R94").Value = 1.85
$ws.Range("S93").Value = 1.7
$ws.Range("S94").Value = 1.95
$ws.Range("T93").Value = 3.25
$ws.Range("T94").Value = 3
$ws.Range("U93").Value = 1.85
$ws.Range("U94").Value = 1.85
$ws.Range("V93").Value = 1.95
$ws.Range("V94").Value = 1.95
$ws.Range("W93").Value = -1
$ws.Range("W94").Value = 0.363
$ws.Range("X93").Value = -1
$ws.Range("X94").Value = -1
$ws.Range("Y93").Value = 0.3
$ws.Range("Y94").Value = -1
$ws.Range("Z93").Value = -1
$ws.Range("Z94").Value = 0.8500000000000001
$ws.Range("AA93").Value = 0.7
$ws.Range("AA94").Value = -1
$ws.Range("AB93").Value = -0.5
$ws.Range("AB94").Value = -1
$ws.Range("AC93").Value = 0.475
$ws.Range("AC94").Value = 0.95

# --- swap rows 108 and 109 ---
$ws.Range("B108").Value = 6100441
$ws.Range("B109").Value = 7265416
$ws.Range("F108").Value = 'Nordic United FC'
$ws.Range("F109").Value = 'Tby FK'
$ws.Range("G108").Value = 'Vasalunds IF'
$ws.Range("G109").Value = 'IF Karlstad Fotboll'
$ws.Range("H108").Value = 1
$ws.Range("H109").Value = 1
$ws.Range("I108").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J108").Value = 'H'
$ws.Range("J109").Value = 'H'
$ws.Range("K108").Value = 2.05
$ws.Range("K109").Value = 2.6
$ws.Range("L108").Value = 3.5
$ws.Range("L109").Value = 3.5
$ws.Range("M108").Value = 3
$ws.Range("M109").Value = 2.25
$ws.Range("N108").Value = 1.95
$ws.Range("N109").Value = 3
$ws.Range("O108").Value = 3.6
$ws.Range("O109").Value = 3.8
$ws.Range("P108").Value = 3.2
$ws.Range("P109").Value = 2.05
$ws.Range("Q108").Value = -0.5
$ws.Range("Q109").Value = 0.25
$ws.Range("R108").Value = 2
$ws.Range("R109").Value = 2
$ws.Range("S108").Value = 1.8
$ws.Range("S109").Value = 1.8
$ws.Range("T108").Value = 3
$ws.Range("T109").Value = 3.25
$ws.Range("U108").Value = 1.925
$ws.Range("U109").Value = 1.95
$ws.Range("V108").Value = 1.875
$ws.Range("V109").Value = 1.85
$ws.Range("W108").Value = 0.95
$ws.Range("W109").Value = 2
$ws.Range("X108").Value = -1
$ws.Range("X109").Value = -1
$ws.Range("Y108").Value = -1
$ws.Range("Y109").Value = -1
$ws.Range("Z108").Value = 1
$ws.Range("Z109").Value = 1
$ws.Range("AA108").Value = -1
$ws.Range("AA109").Value = -1
$ws.Range("AB108").Value = -1
$ws.Range("AB109").Value = -1
$ws.Range("AC108").Value = 0.875
$ws.Range("AC109").Value = 0.8500000000000001

# --- swap rows 117 and 118 ---
$ws.Range("B117").Value = 6100826
$ws.Range("B118").Value = 6100828
$ws.Range("F117").Value = 'Motala AIF FK'
$ws.Range("F118").Value = 'Vasalunds IF'
$ws.Range("G117").Value = 'Nordic United FC'
$ws.Range("G118").Value = 'IF Sylvia'
$ws.Range("H117").Value = 0
$ws.Range("H118").Value = 4
$ws.Range("I117").Value = 2
$ws.Range("I118").Value = 0
$ws.Range("J117").Value = 'A'
$ws.Range("J118").Value = 'H'
$ws.Range("K117").Value = 3.2
$ws.Range("K118").Value = 1.333
$ws.Range("L117").Value = 3.6
$ws.Range("L118").Value = 5
$ws.Range("M117").Value = 1.909
$ws.Range("M118").Value = 6
$ws.Range("N117").Value = 4.5
$ws.Range("N118").Value = 1.142
$ws.Range("O117").Value = 4
$ws.Range("O118").Value = 7
$ws.Range("P117").Value = 1.533
$ws.Range("P118").Value = 12
$ws.Range("Q117").Value = 1
$ws.Range("Q118").Value = -2.25
$ws.Range("R117").Value = 1.875
$ws.Range("R118").Value = 2
$ws.Range("S117").Value = 1.925
$ws.Range("S118").Value = 1.8
$ws.Range("T117").Value = 3.25
$ws.Range("T118").Value = 3.5
$ws.Range("U117").Value = 2
$ws.Range("U118").Value = 1.975
$ws.Range("V117").Value = 1.8
$ws.Range("V118").Value = 1.825
$ws.Range("W117").Value = -1
$ws.Range("W118").Value = 0.1419999999999999
$ws.Range("X117").Value = -1
$ws.Range("X118").Value = -1
$ws.Range("Y117").Value = 0.5329999999999999
$ws.Range("Y118").Value = -1
$ws.Range("Z117").Value = -1
$ws.Range("Z118").Value = 1
$ws.Range("AA117").Value = 0.925
$ws.Range("AA118").Value = -1
$ws.Range("AB117").Value = -1
$ws.Range("AB118").Value = 0.9750000000000001
$ws.Range("AC117").Value = 0.8
$ws.Range("AC118").Value = -1

# ------------------------------------------------------------------
# Part 2: a new fixture row is inserted at row 163 (the former row 163
# shifts down to row 164, keeping its data but with id incremented).
# ------------------------------------------------------------------
$ws.Rows(163).Insert()

# id in column A is (row number - 2); the old row 163 (id 161) is now
# row 164, so its id becomes 162.
$ws.Range("A164").Value = 162

# Give the freshly inserted row 163 the same formatting as the other
# data rows (border/bold/centered id cell, dd date cell).
$ws.Range("A163").Font.Bold = $ws.Range("A162").Font.Bold
$ws.Range("A163").Borders.LineStyle = $ws.Range("A162").Borders.LineStyle
$ws.Range("A163").HorizontalAlignment = $ws.Range("A162").HorizontalAlignment
$ws.Range("A163").VerticalAlignment = $ws.Range("A162").VerticalAlignment
$ws.Range("E163").NumberFormat = $ws.Range("E162").NumberFormat

# Populate the new row 163 with the new fixture data.
$ws.Range("A163").Value = 161
$ws.Range("B163").Value = 7724499
$ws.Range("C163").Value = 'Sweden 1div Norra'
$ws.Range("D163").Value = 'Sweden 1div Norra'
$ws.Range("E163").Value = 45381.375
$ws.Range("F163").Value = 'rebro Syrianska IF'
$ws.Range("G163").Value = 'FC Stockholm Internazionale'
$ws.Range("K163").Value = 3.9
$ws.Range("L163").Value = 3.6
$ws.Range("M163").Value = 1.727
$ws.Range("N163").Value = 3.8
$ws.Range("O163").Value = 3.5
$ws.Range("P163").Value = 1.8
$ws.Range("Q163").Value = 0.5
$ws.Range("R163").Value = 1.975
$ws.Range("S163").Value = 1.825
$ws.Range("T163").Value = 3
$ws.Range("U163").Value = 1.975
$ws.Range("V163").Value = 1.825
$ws.Range("W163").Value = 0
$ws.Range("X163").Value = 0
$ws.Range("Y163").Value = 0
$ws.Range("Z163").Value = 0
$ws.Range("AA163").Value = 0

Write-Output "edit complete"
